{"js": "// Apply the worksheet update: new date header + 25 new division problems,\n// in document order (title paragraph, then each table cell paragraph).\nconst replacements = [\n  \"2024-08-31 Saturday\",\n  \"19\u00f77=2, 5\",\n  \"54\u00f72=27, 0\",\n  \"60\u00f79=6, 6\",\n  \"89\u00f77=12, 5\",\n  \"66\u00f72=33, 0\",\n  \"58\u00f74=14, 2\",\n  \"18\u00f77=2, 4\",\n  \"74\u00f72=37, 0\",\n  \"86\u00f76=14, 2\",\n  \"72\u00f75=14, 2\",\n  \"80\u00f78=10, 0\",\n  \"98\u00f78=12, 2\",\n  \"59\u00f79=6, 5\",\n  \"13\u00f74=3, 1\",\n  \"94\u00f73=31, 1\",\n  \"21\u00f75=4, 1\",\n  \"51\u00f77=7, 2\",\n  \"27\u00f77=3, 6\",\n  \"99\u00f77=14, 1\",\n  \"42\u00f76=7, 0\",\n  \"87\u00f73=29, 0\",\n  \"61\u00f76=10, 1\",\n  \"49\u00f73=16, 1\",\n  \"53\u00f78=6, 5\",\n  \"29\u00f75=5, 4\",\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet idx = 0;\nfor (let i = 0; i < paragraphs.items.length && idx < replacements.length; i++) {\n  const para = paragraphs.items[i];\n  // Every paragraph that originally carried text (the date line and the\n  // 25 populated table cells) gets the next value from the ordered list;\n  // the blank filler-row cells are left untouched.\n  if (para.text !== \"\") {\n    para.insertText(replacements[idx], Word.InsertLocation.replace);\n    idx++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Apply the worksheet update: new date header + 25 new division problems,\n# in document order (title paragraph, then each table cell paragraph).\n$replacements = @(\n    \"2024-08-31 Saturday\",\n    \"19\u00f77=2, 5\",\n    \"54\u00f72=27, 0\",\n    \"60\u00f79=6, 6\",\n    \"89\u00f77=12, 5\",\n    \"66\u00f72=33, 0\",\n    \"58\u00f74=14, 2\",\n    \"18\u00f77=2, 4\",\n    \"74\u00f72=37, 0\",\n    \"86\u00f76=14, 2\",\n    \"72\u00f75=14, 2\",\n    \"80\u00f78=10, 0\",\n    \"98\u00f78=12, 2\",\n    \"59\u00f79=6, 5\",\n    \"13\u00f74=3, 1\",\n    \"94\u00f73=31, 1\",\n    \"21\u00f75=4, 1\",\n    \"51\u00f77=7, 2\",\n    \"27\u00f77=3, 6\",\n    \"99\u00f77=14, 1\",\n    \"42\u00f76=7, 0\",\n    \"87\u00f73=29, 0\",\n    \"61\u00f76=10, 1\",\n    \"49\u00f73=16, 1\",\n    \"53\u00f78=6, 5\",\n    \"29\u00f75=5, 4\"\n)\n\n$d = $word.ActiveDocument\n$count = $d.Paragraphs.Count\n$idx = 0\n\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text\n    # Every real paragraph/cell ends with a paragraph mark (13) and, inside\n    # a table cell, a cell mark (7) as well; strip those off to see whether\n    # there is actual text content.\n    $clean = $t.TrimEnd([char]13, [char]7)\n    if ($clean.Length -gt 0) {\n        $p.Range.Text = $replacements[$idx]\n        $idx = $idx + 1\n    }\n}\n"}
